$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (PREVOST Arthur)
$ws.Range("D6").Value = "231231APR1"
$ws.Range("E6").Value = "NA"
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = 282

# Row 7 (JULIEN Marion)
$ws.Range("D7").Value = "231231MJR1"
$ws.Range("E7").Value = "NA"
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = 132

# Row 12 (PASTIER Jean-pierre)
$ws.Range("D12").Value = "231231JPR1"
$ws.Range("E12").Value = "NA"
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = 54

# Row 23 (TIEDE Alice)
$ws.Range("D23").Value = "231231AT1"
